# "Conditions use for different activity"
#
# 1. Insert a brand-new worksheet named "Sheet1" before "AddAchievement Test Data".
#    It holds a single row of data that used to live (partially) on the
#    "CreateOwn Test Data" sheet (the "Playing Chess" record).
# 2. Extend "CreateOwn Test Data" with new columns (G:O) describing the
#    extra conditions needed to create different kinds of activities
#    (Once / Every Day, time required, date, completeWithin, endAfterDays,
#    expected messages, points, status, toast message, etc.), and rename
#    column F's header to "expectedActivityCreatedMsg".

$wb = $excel.ActiveWorkbook

$createOwn = $wb.Worksheets.Item("CreateOwn Test Data")
$addAchievement = $wb.Worksheets.Item("AddAchievement Test Data")

# --- 1. Insert the new "Sheet1" tab right before "AddAchievement Test Data" ---
$newSheet = $wb.Worksheets.Add($addAchievement)
$newSheet.Name = "Sheet1"

$newSheet.Range("A1").Value = "satish.kshirsagar@gmail.com"
$newSheet.Range("B1").Value = "'1111"
$newSheet.Range("C1").Value = "Playing Chess"
$newSheet.Range("D1").Value = "Playing Chess with friends"
$newSheet.Range("E1").Value = "Chess"
$newSheet.Range("F1").Value = "ACTIVITY CREATED SUCCESSFULLY"

$newSheet.Rows.Item(1).Select() | Out-Null

# --- 2. Extend "CreateOwn Test Data" with the new condition columns ---

# Header row (row 1) - value first ...
$createOwn.Range("F1").Value = "expectedActivityCreatedMsg"
$createOwn.Range("G1").Value = "activityType"
$createOwn.Range("H1").Value = "timeRequired"
$createOwn.Range("I1").Value = "date"
$createOwn.Range("J1").Value = "completeWithin"
$createOwn.Range("K1").Value = "endAfterDays"
$createOwn.Range("L1").Value = "expectedAssignedActivityMsg"
$createOwn.Range("M1").Value = "points"
$createOwn.Range("N1").Value = "statusOfActivity"
$createOwn.Range("O1").Value = "expectedToastMsg"

# ... then apply the same yellow header fill used by the rest of row 1
# (this reuses the existing "header" cell style instead of creating a new one).
$createOwn.Range("G1:O1").Interior.Color = 65535

# Row 2 - "Once" activity data
$createOwn.Range("G2").Value = "Once"
$createOwn.Range("H2").Value = "'30"
$createOwn.Range("I2").Value = "'16"
$createOwn.Range("J2").Value = "DAY"
$createOwn.Range("K2").Value = "'"
$createOwn.Range("L2").Value = "YAY! IT'S ASSIGNED"
$createOwn.Range("M2").Value = "'5"
$createOwn.Range("N2").Value = "Completed"
$createOwn.Range("O2").Value = "Activity Deleted Successfully"

# Row 3 - "Every Day" activity data
$createOwn.Range("G3").Value = "Every Day"
$createOwn.Range("H3").Value = "'60"
$createOwn.Range("I3").Value = "'16"
$createOwn.Range("J3").Value = "'"
$createOwn.Range("K3").Value = "'10"
$createOwn.Range("L3").Value = "YAY! IT'S ASSIGNED"
$createOwn.Range("M3").Value = "'5"
$createOwn.Range("N3").Value = "Completed"
$createOwn.Range("O3").Value = "Activity Deleted Successfully"

# Selection / active sheet bookkeeping to mirror the saved view state
$createOwn.Activate()
$createOwn.Range("J7").Select() | Out-Null

$createActivityError = $wb.Worksheets.Item("Create Activity Error Msg Data")
$createActivityError.Range("F10").Select() | Out-Null

$createOwn.Activate()
